# Updated symbol list on Thu Dec 29 15:38:46 UTC 2022 with GitHub Actions
# Refresh price/volume snapshot data for the crypto ranking sheet.
# Note: Price (column D) values are stored as text in the sheet (not
# numbers), so we prefix numeric-looking values with a leading apostrophe
# to force Excel to keep them as text instead of auto-converting to a
# numeric cell (which would lose trailing zeros / introduce FP noise).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.11"
$ws.Range("D4").Value = "'5.364"
$ws.Range("D5").Value = "'0.05736"
$ws.Range("D6").Value = "'6.473"
$ws.Range("D7").Value = "'3.141"
$ws.Range("D8").Value = "'0.8188"
$ws.Range("D9").Value = "'0.8713"
$ws.Range("D10").Value = "'0.1380"
$ws.Range("D11").Value = "'0.06983"
$ws.Range("D12").Value = "'0.03141"
$ws.Range("D13").Value = "'0.02941"
$ws.Range("D14").Value = "'0.09405"
$ws.Range("D15").Value = "'3.739"
$ws.Range("D16").Value = "'0.001532"
$ws.Range("D17").Value = "'0.04693"
$ws.Range("D18").Value = "'0.0005969"
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("D19").Value = "'0.006236"
$ws.Range("D21").Value = "'0.004795"
$ws.Range("D22").Value = "'0.00008799"
$ws.Range("D27").Value = "'0.1328"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006429"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1058"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002800"
$ws.Range("E43").Value = "42CEJICEJIBestin24h"
$ws.Range("D44").Value = "'0.007518"
$ws.Range("E44").Value = "43LocalTradersLCT"
$ws.Range("D45").Value = "'0.00005264"
$ws.Range("D47").Value = "'0.3998"
$ws.Range("D48").Value = "'0.002486"
$ws.Range("E48").Value = "47BOLOBOLO"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("D50").Value = "'0.0002000"
